$d = $word.ActiveDocument

# Locate the paragraph that ends the section to keep: "...MÁQUINASEscola PRO-TEC"
$anchor = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*PRO-TEC*") {
        $anchor = $p
    }
}

# Locate the last paragraph of the trailing block that must be removed (the
# copyright / footer notice paragraph).
$lastToRemove = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Creative Commons Attribution*") {
        $lastToRemove = $p
    }
}

# Remove everything between the end of the anchor paragraph and the end of
# the last paragraph to remove (inclusive of their paragraph marks). This
# deletes the blank paragraph, the "Ver no Jupiter..." paragraph, and the
# "© 2020 ..." paragraph in one shot, while leaving the paragraph mark of
# the anchor paragraph (and everything after the removed block) untouched.
$startPos = $anchor.Range.End
$endPos = $lastToRemove.Range.End
$d.Range($startPos, $endPos).Delete()
